$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (the "Förändrad" date column) for rows 2 through 17
# from Excel date serial 45181 (2023-09-12) to 45182 (2023-09-13).
$ws.Range("C2:C17").Value = 45182
